# Commit: "adding year from admin Type: SAVE."
# Rolls the FICA rates table forward by one year: each existing year-group
# shifts down to the next row-group (2025->2026, 2024->2025, 2023->2024,
# 2022->2023, 2021->2022), and a new trailing 2021 group is appended, with a
# new blank spacer row after it. Also refreshes the header formula-name
# label in B2 and normalizes D19 to a real number.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Configs")

# Header label (B2) - updated function signature text
$ws.Range("B2").Value = "SmartRules PaymentToolsFicaDetails FicaTaxRates(PaymentToolsFicaInput input)"

# --- Year 2026 (was 2025 block: rows 5-7) ---
$ws.Range("B5").Value = 2026
$ws.Range("B6").Value = 2026
$ws.Range("B7").Value = 2026

# --- Year 2025 (was 2024 block: rows 8-10) ---
$ws.Range("B8").Value = 2025
$ws.Range("B9").Value = 2025
$ws.Range("B10").Value = 2025

# --- Year 2024 (was 2023 block: rows 11-13) ---
$ws.Range("B11").Value = 2024
$ws.Range("B12").Value = 2024
$ws.Range("B13").Value = 2024
$ws.Range("E13").Value = 168600

# --- Year 2023 (was 2022 block: rows 14-16) ---
$ws.Range("B14").Value = 2023
$ws.Range("B15").Value = 2023
$ws.Range("B16").Value = 2023
$ws.Range("E16").Value = 152000

# --- Year 2022 (was 2021 block: rows 17-19) ---
$ws.Range("B17").Value = 2022
$ws.Range("B18").Value = 2022
$ws.Range("B19").Value = 2022
$ws.Range("D19").Value = 6.2
$ws.Range("E19").Value = 147000

# --- New trailing year 2021 block (rows 20-22), appended at bottom ---
# Clear old leftover "template" formatting from the prior blank row 20 so
# the new cells pick up the plain default column style, matching a brand
# new row.
$ws.Range("B20:E22").Clear()

$ws.Range("B20").Value = 2021
$ws.Range("C20").Value = "AFMT"
$ws.Range("D20").Value = 0.9
$ws.Range("E20").Value = 200000

$ws.Range("B21").Value = 2021
$ws.Range("C21").Value = "FMT"
$ws.Range("D21").Value = 1.45

$ws.Range("B22").Value = 2021
$ws.Range("C22").Value = "FSST"
$ws.Range("D22").Value = 6.2
$ws.Range("E22").Value = 142800

# --- New trailing blank spacer row 23 ---
$ws.Range("B23:E23").Clear()
